$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 16.71895933333333
$ws.Range("H2").Value = 50.156878
$ws.Range("I2").Value = 0.02912144738161902
$ws.Range("J2").Value = 0.03059269312988411
$ws.Range("M2").Value = 6.712486666666666
$ws.Range("N2").Value = 20.13746
$ws.Range("O2").Value = 0.6330487633990675
$ws.Range("P2").Value = 0.6414503882251803
$ws.Range("Q2").Value = 112.2257916055422
$ws.Range("R2").Value = 1010.03212444988
$ws.Range("S2").Value = 0.01843529625332493
$ws.Range("T2").Value = 0.01962369488501797
$ws.Range("G3").Value = 16.71895933333333
$ws.Range("H3").Value = 50.156878
$ws.Range("I3").Value = 0.02912144738161902
$ws.Range("J3").Value = 0.03059269312988411
$ws.Range("O3").Value = 0.290741083484562
$ws.Range("P3").Value = 0.2945997080427384
$ws.Range("Q3").Value = 51.54207721868534
$ws.Range("R3").Value = 463.878694968168
$ws.Range("S3").Value = 0.008466801164370572
$ws.Range("T3").Value = 0.009012598464304947
$ws.Range("G4").Value = 16.71895933333333
$ws.Range("H4").Value = 50.156878
$ws.Range("I4").Value = 0.02912144738161902
$ws.Range("J4").Value = 0.03059269312988411
$ws.Range("M4").Value = 0.2495096666666667
$ws.Range("N4").Value = 0.748529
$ws.Range("O4").Value = 0.02353103905946135
$ws.Range("P4").Value = 0.02384333563656022
$ws.Range("Q4").Value = 4.171541970273555
$ws.Range("R4").Value = 37.543877732462
$ws.Range("S4").Value = 0.0006852579158049256
$ws.Range("T4").Value = 0.000729431850322117
$ws.Range("G5").Value = 16.71895933333333
$ws.Range("H5").Value = 50.156878
$ws.Range("I5").Value = 0.02912144738161902
$ws.Range("J5").Value = 0.03059269312988411
$ws.Range("M5").Value = 0.4166465
$ws.Range("N5").Value = 0.8332930000000001
$ws.Range("O5").Value = 0.03929356804674715
$ws.Range("P5").Value = 0.02654337331298611
$ws.Range("Q5").Value = 6.965895889875667
$ws.Range("R5").Value = 41.795375339254
$ws.Range("S5").Value = 0.001144285574309413
$ws.Range("T5").Value = 0.0008120332743961395
$ws.Range("G6").Value = 16.71895933333333
$ws.Range("H6").Value = 50.156878
$ws.Range("I6").Value = 0.02912144738161902
$ws.Range("J6").Value = 0.03059269312988411
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.1419326666666667
$ws.Range("N6").Value = 0.425798
$ws.Range("O6").Value = 0.01338554601016197
$ws.Range("P6").Value = 0.01356319478253491
$ws.Range("Q6").Value = 2.372966482071556
$ws.Range("R6").Value = 21.356698338644
$ws.Range("S6").Value = 0.000389806473809172
$ws.Range("T6").Value = 0.0004149346558429356
$ws.Range("I7").Value = 0.2708539632042961
$ws.Range("J7").Value = 0.2845377865576845
$ws.Range("M7").Value = 6.712486666666666
$ws.Range("N7").Value = 20.13746
$ws.Range("O7").Value = 0.6330487633990675
$ws.Range("P7").Value = 0.6414503882251803
$ws.Range("Q7").Value = 1043.79428782398
$ws.Range("R7").Value = 9394.148590415818
$ws.Range("S7").Value = 0.1714637664682161
$ws.Range("T7").Value = 0.1825168736521602
$ws.Range("I8").Value = 0.2708539632042961
$ws.Range("J8").Value = 0.2845377865576845
$ws.Range("O8").Value = 0.290741083484562
$ws.Range("P8").Value = 0.2945997080427384
$ws.Range("S8").Value = 0.07874837472810471
$ws.Range("T8").Value = 0.08382474884702087
$ws.Range("I9").Value = 0.2708539632042961
$ws.Range("J9").Value = 0.2845377865576845
$ws.Range("M9").Value = 0.2495096666666667
$ws.Range("N9").Value = 0.748529
$ws.Range("O9").Value = 0.02353103905946135
$ws.Range("P9").Value = 0.02384333563656022
$ws.Range("Q9").Value = 38.798850225927
$ws.Range("R9").Value = 349.189652033343
$ws.Range("S9").Value = 0.006373475187570198
$ws.Range("T9").Value = 0.006784329946178807
$ws.Range("I10").Value = 0.2708539632042961
$ws.Range("J10").Value = 0.2845377865576845
$ws.Range("M10").Value = 0.4166465
$ws.Range("N10").Value = 0.8332930000000001
$ws.Range("O10").Value = 0.03929356804674715
$ws.Range("P10").Value = 0.02654337331298611
$ws.Range("Q10").Value = 64.78869282548851
$ws.Range("R10").Value = 388.7321569529311
$ws.Range("S10").Value = 0.01064281863389916
$ws.Range("T10").Value = 0.007552592690251382
$ws.Range("I11").Value = 0.2708539632042961
$ws.Range("J11").Value = 0.2845377865576845
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.1419326666666667
$ws.Range("N11").Value = 0.425798
$ws.Range("O11").Value = 0.01338554601016197
$ws.Range("P11").Value = 0.01356319478253491
$ws.Range("Q11").Value = 22.070584878474
$ws.Range("R11").Value = 198.635263906266
$ws.Range("S11").Value = 0.003625528186505821
$ws.Range("T11").Value = 0.003859241422073217
$ws.Range("G12").Value = 194.8548433333333
$ws.Range("H12").Value = 584.56453
$ws.Range("I12").Value = 0.3394024086099587
$ws.Range("J12").Value = 0.3565493705749576
$ws.Range("M12").Value = 6.712486666666666
$ws.Range("N12").Value = 20.13746
$ws.Range("O12").Value = 0.6330487633990675
$ws.Range("P12").Value = 0.6414503882251803
$ws.Range("Q12").Value = 1307.960537810422
$ws.Range("R12").Value = 11771.6448402938
$ws.Range("S12").Value = 0.2148582750651994
$ws.Range("T12").Value = 0.2287087321767502
$ws.Range("G13").Value = 194.8548433333333
$ws.Range("H13").Value = 584.56453
$ws.Range("I13").Value = 0.3394024086099587
$ws.Range("J13").Value = 0.3565493705749576
$ws.Range("O13").Value = 0.290741083484562
$ws.Range("P13").Value = 0.2945997080427384
$ws.Range("Q13").Value = 600.7086434798533
$ws.Range("R13").Value = 5406.37779131868
$ws.Range("S13").Value = 0.09867822401652943
$ws.Range("T13").Value = 0.1050393404742046
$ws.Range("G14").Value = 194.8548433333333
$ws.Range("H14").Value = 584.56453
$ws.Range("I14").Value = 0.3394024086099587
$ws.Range("J14").Value = 0.3565493705749576
$ws.Range("M14").Value = 0.2495096666666667
$ws.Range("N14").Value = 0.748529
$ws.Range("O14").Value = 0.02353103905946135
$ws.Range("P14").Value = 0.02384333563656022
$ws.Range("Q14").Value = 48.61816700848555
$ws.Range("R14").Value = 437.56350307637
$ws.Range("S14").Value = 0.007986491333876201
$ws.Range("T14").Value = 0.008501326313623003
$ws.Range("G15").Value = 194.8548433333333
$ws.Range("H15").Value = 584.56453
$ws.Range("I15").Value = 0.3394024086099587
$ws.Range("J15").Value = 0.3565493705749576
$ws.Range("M15").Value = 0.4166465
$ws.Range("N15").Value = 0.8332930000000001
$ws.Range("O15").Value = 0.03929356804674715
$ws.Range("P15").Value = 0.02654337331298611
$ws.Range("Q15").Value = 81.18558848288167
$ws.Range("R15").Value = 487.11353089729
$ws.Range("S15").Value = 0.01333633163794529
$ws.Range("T15").Value = 0.009464023047681324
$ws.Range("G16").Value = 194.8548433333333
$ws.Range("H16").Value = 584.56453
$ws.Range("I16").Value = 0.3394024086099587
$ws.Range("J16").Value = 0.3565493705749576
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.1419326666666667
$ws.Range("N16").Value = 0.425798
$ws.Range("O16").Value = 0.01338554601016197
$ws.Range("P16").Value = 0.01356319478253491
$ws.Range("Q16").Value = 27.65626752721555
$ws.Range("R16").Value = 248.90640774494
$ws.Range("S16").Value = 0.004543086556408395
$ws.Range("T16").Value = 0.004835948562698369
$ws.Range("G17").Value = 82.82950199999999
$ws.Range("H17").Value = 165.659004
$ws.Range("I17").Value = 0.1442742299952585
$ws.Range("J17").Value = 0.1010420758958371
$ws.Range("M17").Value = 6.712486666666666
$ws.Range("N17").Value = 20.13746
$ws.Range("O17").Value = 0.6330487633990675
$ws.Range("P17").Value = 0.6414503882251803
$ws.Range("Q17").Value = 555.9919277816399
$ws.Range("R17").Value = 3335.951566689839
$ws.Range("S17").Value = 0.09133262288885104
$ws.Range("T17").Value = 0.06481347881046286
$ws.Range("G18").Value = 82.82950199999999
$ws.Range("H18").Value = 165.659004
$ws.Range("I18").Value = 0.1442742299952585
$ws.Range("J18").Value = 0.1010420758958371
$ws.Range("O18").Value = 0.290741083484562
$ws.Range("P18").Value = 0.2945997080427384
$ws.Range("Q18").Value = 255.351095899704
$ws.Range("R18").Value = 1532.106575398224
$ws.Range("S18").Value = 0.04194644594772234
$ws.Range("T18").Value = 0.02976696605894583
$ws.Range("G19").Value = 82.82950199999999
$ws.Range("H19").Value = 165.659004
$ws.Range("I19").Value = 0.1442742299952585
$ws.Range("J19").Value = 0.1010420758958371
$ws.Range("M19").Value = 0.2495096666666667
$ws.Range("N19").Value = 0.748529
$ws.Range("O19").Value = 0.02353103905946135
$ws.Range("P19").Value = 0.02384333563656022
$ws.Range("Q19").Value = 20.666761434186
$ws.Range("R19").Value = 124.000568605116
$ws.Range("S19").Value = 0.003394922541292138
$ws.Range("T19").Value = 0.002409180128999237
$ws.Range("G20").Value = 82.82950199999999
$ws.Range("H20").Value = 165.659004
$ws.Range("I20").Value = 0.1442742299952585
$ws.Range("J20").Value = 0.1010420758958371
$ws.Range("M20").Value = 0.4166465
$ws.Range("N20").Value = 0.8332930000000001
$ws.Range("O20").Value = 0.03929356804674715
$ws.Range("P20").Value = 0.02654337331298611
$ws.Range("Q20").Value = 34.510622105043
$ws.Range("R20").Value = 138.042488420172
$ws.Range("S20").Value = 0.005669049273710738
$ws.Range("T20").Value = 0.002681997540822281
$ws.Range("G21").Value = 82.82950199999999
$ws.Range("H21").Value = 165.659004
$ws.Range("I21").Value = 0.1442742299952585
$ws.Range("J21").Value = 0.1010420758958371
$ws.Range("K21").Value = 2
$ws.Range("L21").Value = 0.6666666666666666
$ws.Range("M21").Value = 0.1419326666666667
$ws.Range("N21").Value = 0.425798
$ws.Range("O21").Value = 0.01338554601016197
$ws.Range("P21").Value = 0.01356319478253491
$ws.Range("Q21").Value = 11.756212097532
$ws.Range("R21").Value = 70.537272585192
$ws.Range("S21").Value = 0.001931189343682222
$ws.Range("T21").Value = 0.001370453356606914
$ws.Range("G22").Value = 124.2078576666667
$ws.Range("H22").Value = 372.623573
$ws.Range("I22").Value = 0.2163479508088675
$ws.Range("J22").Value = 0.2272780738416368
$ws.Range("M22").Value = 6.712486666666666
$ws.Range("N22").Value = 20.13746
$ws.Range("O22").Value = 0.6330487633990675
$ws.Range("P22").Value = 0.6414503882251803
$ws.Range("Q22").Value = 833.743588482731
$ws.Range("R22").Value = 7503.692296344579
$ws.Range("S22").Value = 0.1369588027234759
$ws.Range("T22").Value = 0.1457876087007891
$ws.Range("G23").Value = 124.2078576666667
$ws.Range("H23").Value = 372.623573
$ws.Range("I23").Value = 0.2163479508088675
$ws.Range("J23").Value = 0.2272780738416368
$ws.Range("O23").Value = 0.290741083484562
$ws.Range("P23").Value = 0.2945997080427384
$ws.Range("Q23").Value = 382.9144424233987
$ws.Range("R23").Value = 3446.229981810588
$ws.Range("S23").Value = 0.06290123762783487
$ws.Range("T23").Value = 0.06695605419826214
$ws.Range("G24").Value = 124.2078576666667
$ws.Range("H24").Value = 372.623573
$ws.Range("I24").Value = 0.2163479508088675
$ws.Range("J24").Value = 0.2272780738416368
$ws.Range("M24").Value = 0.2495096666666667
$ws.Range("N24").Value = 0.748529
$ws.Range("O24").Value = 0.02353103905946135
$ws.Range("P24").Value = 0.02384333563656022
$ws.Range("Q24").Value = 30.99106116379077
$ws.Range("R24").Value = 278.919550474117
$ws.Range("S24").Value = 0.005090892080917886
$ws.Range("T24").Value = 0.005419067397437066
$ws.Range("G25").Value = 124.2078576666667
$ws.Range("H25").Value = 372.623573
$ws.Range("I25").Value = 0.2163479508088675
$ws.Range("J25").Value = 0.2272780738416368
$ws.Range("M25").Value = 0.4166465
$ws.Range("N25").Value = 0.8332930000000001
$ws.Range("O25").Value = 0.03929356804674715
$ws.Range("P25").Value = 0.02654337331298611
$ws.Range("Q25").Value = 51.75076916931483
$ws.Range("R25").Value = 310.504615015889
$ws.Range("S25").Value = 0.008501082926882541
$ws.Range("T25").Value = 0.00603272675983499
$ws.Range("G26").Value = 124.2078576666667
$ws.Range("H26").Value = 372.623573
$ws.Range("I26").Value = 0.2163479508088675
$ws.Range("J26").Value = 0.2272780738416368
$ws.Range("K26").Value = 2
$ws.Range("L26").Value = 0.6666666666666666
$ws.Range("M26").Value = 0.1419326666666667
$ws.Range("N26").Value = 0.425798
$ws.Range("O26").Value = 0.01338554601016197
$ws.Range("P26").Value = 0.01356319478253491
$ws.Range("Q26").Value = 17.62915245958378
$ws.Range("R26").Value = 158.662372136254
$ws.Range("S26").Value = 0.002895935449756354
$ws.Range("T26").Value = 0.003082616785313472
